# Refresh the stock data in CompanyStockInfo:
#  - Abbott India (row 2): small EUR column refresh
#  - Rows 3-6 reshuffled with Just Dial / DLF / Aditya Birla F / Ceat,
#    where DLF is a brand new entry replacing the old "DFL / Company Not Listed"
#    placeholder row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Abbott India: refreshed EUR quotes
$ws.Range("F2").Value2 = 319.78
$ws.Range("G2").Value2 = 310.35
$ws.Range("H2").Value2 = 310.61

# Row 3 - Just Dial
$ws.Range("A3").Value2 = "Just Dial"
$ws.Range("B3").Value2 = 1199.6
$ws.Range("C3").Value2 = "18/10/2024 15:59"
$ws.Range("D3").Value2 = 1194.95
$ws.Range("E3").Value2 = 1194.4
$ws.Range("F3").Value2 = 13.13
$ws.Range("G3").Value2 = 13.08
$ws.Range("H3").Value2 = 13.07

# Row 4 - DLF (new entry, replaces the old "DFL / Company Not Listed" row)
$ws.Range("A4").Value2 = "DLF"
$ws.Range("B4").Value2 = 875.15
$ws.Range("C4").Value2 = "18/10/2024 15:59"
$ws.Range("D4").Value2 = 860.9
$ws.Range("E4").Value2 = 861
$ws.Range("F4").Value2 = 9.58
$ws.Range("G4").Value2 = 9.42
$ws.Range("H4").Value2 = 9.42

# Row 5 - Aditya Birla F
$ws.Range("A5").Value2 = "Aditya Birla F"
$ws.Range("B5").Value2 = 334.05
$ws.Range("C5").Value2 = "18/10/2024 15:59"
$ws.Range("D5").Value2 = 330.9
$ws.Range("E5").Value2 = 330.7
$ws.Range("F5").Value2 = 3.65
$ws.Range("G5").Value2 = 3.62
$ws.Range("H5").Value2 = 3.62

# Row 6 - Ceat
$ws.Range("A6").Value2 = "Ceat"
$ws.Range("B6").Value2 = 2972.85
$ws.Range("C6").Value2 = "18/10/2024 16:00"
$ws.Range("D6").Value2 = 2799
$ws.Range("E6").Value2 = 2891.75
$ws.Range("F6").Value2 = 32.54
$ws.Range("G6").Value2 = 30.64
$ws.Range("H6").Value2 = 31.65
